# Apply the "update text and figures" edit to slide 1 of the presentation.
#
# Changes:
#   1. Group 10 > Pentagon 3 (id 15): "Situational Report" -> "Situation Update"
#   2. Group 10 > Rectangle 18 (id 19): first line "CRF" -> "CFR"
#   3. Rectangle 1 (id 2, title): fix typo + add "2" to the report number,
#      splitting the sentence into two runs, and widen the title textbox
#      so the longer text still fits (autosized width: 9208868 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Situational Report" -> "Situation Update" (inside Group 10) ---
$group10 = $s.Shapes.Item(2)
$pentagon = $group10.GroupItems.Item(2)
$pentagon.TextFrame.TextRange.Runs(1, 1).Text = "Situation Update"

# --- 2. "CRF" -> "CFR" (inside the same group) ---
$rect18 = $group10.GroupItems.Item(6)
$rect18.TextFrame.TextRange.Runs(1, 1).Text = "CFR"

# --- 3. Title text fix + split into two runs, then widen the textbox ---
$title = $s.Shapes.Item(3)
$titleRange = $title.TextFrame.TextRange
$firstRun = $titleRange.Runs(1, 1)
$firstRun.Text = "SITUATIONAL REPORT 2: "
$firstRun.InsertAfter("PNEUMONIC PLAGUE OUTBREAK IN NORTHEAST INDIA") | Out-Null

# Widen the (auto-fit) title shape to match the new, longer text.
$title.Width = 9208868 / 12700
